# Marksheet update: fix scoring summary, populate "Student Ans" column (A)
# with correct/incorrect answers, and drop the redundant duplicate
# "Student Ans/Correct Ans" blocks that used to live in columns D:E (for
# most rows) and G:H (for all rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header cells in column A (rows 10-12) gain the "mtitleStyle" look
#    (matches style already used by the row-9 header cells).
# ---------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 2. Updated summary numbers.
# ---------------------------------------------------------------------
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 18
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 36
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "35/112"

# ---------------------------------------------------------------------
# 3. Drop the duplicate "Student Ans"/"Correct Ans" block in columns G:H
#    entirely (all rows 15-21).
# ---------------------------------------------------------------------
$ws.Range("G15:H21").Clear()

# ---------------------------------------------------------------------
# 4. Drop the duplicate block in columns D:E for rows 19-40 (rows 16-18
#    keep their existing D:E content).
# ---------------------------------------------------------------------
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------
# 5. Row 17 is special: instead of clearing D17, it gets repurposed to
#    show the correct answer ("Option C") highlighted with the
#    "correctStyle" (green) look.
# ---------------------------------------------------------------------
$ws.Range("B10").Copy()
$ws.Range("D17").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D17").Value = "Option C"

# ---------------------------------------------------------------------
# 6. Populate column A ("Student Ans") for the rows whose answer is
#    known, styled green (correctStyle) when it matches the correct
#    answer in column B, or red (incorrectStyle) otherwise.
# ---------------------------------------------------------------------
$correctRows = @(22, 27, 29, 30, 32, 35, 38, 40)
foreach ($r in $correctRows) {
    $ws.Range("B10").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("A$r").Value = $ws.Range("B$r").Value2
}

$incorrectRows = @(24)
foreach ($r in $incorrectRows) {
    $ws.Range("C10").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("A$r").Value = "Option D"
}
